$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 295.25
$ws.Range("I5").Value = 172.85715
$ws.Range("J5").Value = 466.6
$ws.Range("K5").Value = 172.85715
$ws.Range("L5").Value = 466.6
$ws.Range("M5").Value = -57.85714999999999
$ws.Range("N5").Value = -696.6
$ws.Range("H28").Value = 552.8148
$ws.Range("I28").Value = 538
$ws.Range("J28").Value = 588
$ws.Range("K28").Value = 538
$ws.Range("L28").Value = 588
$ws.Range("M28").Value = -53
$ws.Range("N28").Value = -1558
$ws.Range("H116").Value = 675913.2
$ws.Range("I116").Value = 1432270.2
$ws.Range("J116").Value = 14100.75
$ws.Range("K116").Value = 1432270.2
$ws.Range("L116").Value = 14100.75
$ws.Range("M116").Value = -1428828.2
$ws.Range("N116").Value = -20984.75
$ws.Range("H132").Value = 231019.94
$ws.Range("I132").Value = 3939.2646
$ws.Range("J132").Value = 1003094.2
$ws.Range("K132").Value = 11817.7938
$ws.Range("L132").Value = 3009282.6
$ws.Range("M132").Value = -9287.793799999999
$ws.Range("N132").Value = -3014342.6
$ws.Range("H133").Value = 47806
$ws.Range("J133").Value = 47806
$ws.Range("L133").Value = 47806
$ws.Range("N133").Value = -57926
$ws.Range("H137").Value = 3528.76
$ws.Range("I137").Value = 2337.6
$ws.Range("K137").Value = 7012.799999999999
$ws.Range("M137").Value = -4462.799999999999
$ws.Range("H141").Value = 6013.1626
$ws.Range("I141").Value = 6280.421
$ws.Range("J141").Value = 3982
$ws.Range("K141").Value = 18841.263
$ws.Range("L141").Value = 11946
$ws.Range("M141").Value = -13661.263
$ws.Range("N141").Value = -22306

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 17091
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 17091
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 17091
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -17431
$ws.Range("H32").Value = 5411.0894
$ws.Range("I32").Value = 4287.8335
$ws.Range("K32").Value = 4287.8335
$ws.Range("M32").Value = -4000.8335
$ws.Range("H45").Value = 1455.5264
$ws.Range("I45").Value = 1405.6666
$ws.Range("J45").Value = 1541
$ws.Range("K45").Value = 1405.6666
$ws.Range("L45").Value = 1541
$ws.Range("M45").Value = -1028.6666
$ws.Range("N45").Value = -2295
$ws.Range("H61").Value = 1034.409
$ws.Range("I61").Value = 703
$ws.Range("K61").Value = 703
$ws.Range("M61").Value = -491
$ws.Range("H74").Value = 4043.2144
$ws.Range("I74").Value = 4052.6086
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 4052.6086
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -3178.6086
$ws.Range("N74").Value = -5748
$ws.Range("H76").Value = 40000
$ws.Range("J76").Value = 40000
$ws.Range("L76").Value = 40000
$ws.Range("N76").Value = -40676
$ws.Range("H77").Value = 4043.2144
$ws.Range("I77").Value = 4052.6086
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 20263.043
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -15895.043
$ws.Range("N77").Value = -28736
$ws.Range("H79").Value = 40000
$ws.Range("J79").Value = 40000
$ws.Range("L79").Value = 40000
$ws.Range("N79").Value = -42340
$ws.Range("H122").Value = 2214.1
$ws.Range("I122").Value = 1384.4286
$ws.Range("J122").Value = 4150
$ws.Range("K122").Value = 4153.2858
$ws.Range("L122").Value = 12450
$ws.Range("M122").Value = -1703.2858
$ws.Range("N122").Value = -17350
$ws.Range("H132").Value = 1927.5588
$ws.Range("I132").Value = 824.61536
$ws.Range("J132").Value = 5512.125
$ws.Range("K132").Value = 2473.84608
$ws.Range("L132").Value = 16536.375
$ws.Range("M132").Value = 56.15391999999974
$ws.Range("N132").Value = -21596.375
$ws.Range("H136").Value = 1034.409
$ws.Range("I136").Value = 703
$ws.Range("K136").Value = 2109
$ws.Range("M136").Value = 441

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 359.625
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 359.625
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 359.625
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -809.625
$ws.Range("H67").Value = 359.625
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 359.625
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 359.625
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -1919.625
$ws.Range("H86").Value = 1971.619
$ws.Range("I86").Value = 1414.2858
$ws.Range("K86").Value = 1414.2858
$ws.Range("M86").Value = -291.2858000000001
$ws.Range("H89").Value = 1971.619
$ws.Range("I89").Value = 1414.2858
$ws.Range("K89").Value = 7071.429
$ws.Range("M89").Value = -1455.429
$ws.Range("H134").Value = 1336.8448
$ws.Range("I134").Value = 914.125
$ws.Range("J134").Value = 3365.9
$ws.Range("K134").Value = 2742.375
$ws.Range("L134").Value = 10097.7
$ws.Range("M134").Value = -207.375
$ws.Range("N134").Value = -15167.7
$ws.Range("H135").Value = 43780
$ws.Range("J135").Value = 43780
$ws.Range("L135").Value = 43780
$ws.Range("N135").Value = -53920

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21742218
$ws.Range("I31").Value = 1275.5385
$ws.Range("J31").Value = 50005440
$ws.Range("K31").Value = 1275.5385
$ws.Range("L31").Value = 50005440
$ws.Range("M31").Value = -980.5385000000001
$ws.Range("N31").Value = -50006030
$ws.Range("H34").Value = 21742218
$ws.Range("I34").Value = 1275.5385
$ws.Range("J34").Value = 50005440
$ws.Range("K34").Value = 1275.5385
$ws.Range("L34").Value = 50005440
$ws.Range("M34").Value = -1073.5385
$ws.Range("N34").Value = -50005844
$ws.Range("H58").Value = 1895.5343
$ws.Range("I58").Value = 1652.0793
$ws.Range("K58").Value = 1652.0793
$ws.Range("M58").Value = -1449.0793
$ws.Range("H132").Value = 2427.7856
$ws.Range("I132").Value = 1824.1459
$ws.Range("K132").Value = 5472.4377
$ws.Range("M132").Value = -2942.4377
$ws.Range("H134").Value = 4445.6875
$ws.Range("I134").Value = 5191.773
$ws.Range("J134").Value = 2804.3
$ws.Range("K134").Value = 15575.319
$ws.Range("L134").Value = 8412.900000000001
$ws.Range("M134").Value = -13040.319
$ws.Range("N134").Value = -13482.9
$ws.Range("H136").Value = 1895.5343
$ws.Range("I136").Value = 1652.0793
$ws.Range("K136").Value = 4956.2379
$ws.Range("M136").Value = -2406.2379

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 790.3333
$ws.Range("J113").Value = 1080.8
$ws.Range("L113").Value = 3242.4
$ws.Range("N113").Value = -7582.4
$ws.Range("H133").Value = 2593
$ws.Range("I133").Value = 2632.8572
$ws.Range("K133").Value = 7898.571599999999
$ws.Range("M133").Value = -2838.571599999999
$ws.Range("H134").Value = 4796.913
$ws.Range("I134").Value = 4682.9
$ws.Range("J134").Value = 4884.615
$ws.Range("K134").Value = 14048.7
$ws.Range("L134").Value = 14653.845
$ws.Range("M134").Value = -8978.699999999999
$ws.Range("N134").Value = -24793.845

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4936074
$ws.Range("I11").Value = 9333333
$ws.Range("J11").Value = 2004567.8
$ws.Range("K11").Value = 9333333
$ws.Range("L11").Value = 2004567.8
$ws.Range("M11").Value = -9333194
$ws.Range("N11").Value = -2004845.8
$ws.Range("H12").Value = 14888.5
$ws.Range("I12").Value = 7777
$ws.Range("K12").Value = 7777
$ws.Range("M12").Value = -7637
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50518
$ws.Range("H132").Value = 2984.318
$ws.Range("I132").Value = 1609.5333
$ws.Range("J132").Value = 5930.2856
$ws.Range("K132").Value = 4828.5999
$ws.Range("L132").Value = 17790.8568
$ws.Range("M132").Value = -2298.5999
$ws.Range("N132").Value = -22850.8568
$ws.Range("H133").Value = 38527.777
$ws.Range("J133").Value = 38527.777
$ws.Range("L133").Value = 38527.777
$ws.Range("N133").Value = -48647.777

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 7999.5
$ws.Range("J25").Value = 7999.5
$ws.Range("L25").Value = 7999.5
$ws.Range("N25").Value = -8459.5
$ws.Range("H136").Value = 1775.2128
$ws.Range("I136").Value = 1028.125
$ws.Range("J136").Value = 6044.2856
$ws.Range("K136").Value = 3084.375
$ws.Range("L136").Value = 18132.8568
$ws.Range("M136").Value = -534.375
$ws.Range("N136").Value = -23232.8568

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7093856
$ws.Range("I132").Value = 930.4545000000001
$ws.Range("K132").Value = 2791.3635
$ws.Range("M132").Value = -261.3635000000004
$ws.Range("H136").Value = 2381.0789
$ws.Range("I136").Value = 713.5172
$ws.Range("J136").Value = 7754.3335
$ws.Range("K136").Value = 2140.5516
$ws.Range("L136").Value = 23263.0005
$ws.Range("M136").Value = 409.4484000000002
$ws.Range("N136").Value = -28363.0005
